$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "0 поворотов, расстояние концов 6181 м"
$ws.Range("D3").Value = "0 поворотов, расстояние концов 6181 м"
$ws.Range("D4").Value = "3 поворотов, расстояние концов 18895 м"
$ws.Range("D5").Value = "3 поворотов, расстояние концов 18895 м"
$ws.Range("D6").Value = "1 поворотов, расстояние концов 11450 м"
$ws.Range("D7").Value = "2 поворотов, расстояние концов 11450 м"
$ws.Range("D8").Value = "1 поворотов, расстояние концов 18895 м"
$ws.Range("D9").Value = "1 поворотов, расстояние концов 18895 м"
$ws.Range("D10").Value = "0 поворотов, расстояние концов 13747 м"
$ws.Range("D11").Value = "1 поворотов, расстояние концов 13278 м"
$ws.Range("D12").Value = "0 поворотов, расстояние концов 15372 м"
$ws.Range("D13").Value = "0 поворотов, расстояние концов 15196 м"
$ws.Range("D14").Value = "0 поворотов, расстояние концов 5583 м"
$ws.Range("D15").Value = "0 поворотов, расстояние концов 5560 м"
$ws.Range("D16").Value = "2 поворотов, расстояние концов 11915 м"
$ws.Range("D17").Value = "0 поворотов, расстояние концов 12284 м"
$ws.Range("D18").Value = "2 поворотов, расстояние концов 12500 м"
$ws.Range("D19").Value = "2 поворотов, расстояние концов 13085 м"
$ws.Range("D20").Value = "8 острых углов < 60°"
$ws.Range("D21").Value = "6 острых углов < 60°"
$ws.Range("D22").Value = "8 острых углов < 60°"
$ws.Range("D23").Value = "1 поворотов, расстояние концов 10814 м"
$ws.Range("D24").Value = "2 поворотов, расстояние концов 7369 м"
$ws.Range("D25").Value = "1 поворотов, расстояние концов 7362 м"
$ws.Range("D26").Value = "0 поворотов, расстояние концов 10679 м"
$ws.Range("D27").Value = "0 поворотов, расстояние концов 10702 м"
$ws.Range("D28").Value = "0 поворотов, расстояние концов 14019 м"
$ws.Range("D29").Value = "2 поворотов, расстояние концов 14019 м"
$ws.Range("D30").Value = "0 поворотов, расстояние концов 12493 м"
$ws.Range("D31").Value = "0 поворотов, расстояние концов 12493 м"
$ws.Range("D32").Value = "2 поворотов, расстояние концов 14834 м"
$ws.Range("D33").Value = "1 поворотов, расстояние концов 14728 м"
$ws.Range("D34").Value = "0 поворотов, расстояние концов 10847 м"
$ws.Range("D35").Value = "2 поворотов, расстояние концов 10847 м"
$ws.Range("D36").Value = "0 поворотов, расстояние концов 12656 м"
$ws.Range("D37").Value = "2 поворотов, расстояние концов 12656 м"
$ws.Range("D38").Value = "3 поворотов, расстояние концов 13291 м"
$ws.Range("D39").Value = "1 поворотов, расстояние концов 13291 м"
$ws.Range("D40").Value = "1 поворотов, расстояние концов 13291 м"
$ws.Range("D41").Value = "5 острых углов < 60°"
$ws.Range("D42").Value = "0 поворотов, расстояние концов 6995 м"
$ws.Range("D43").Value = "0 поворотов, расстояние концов 6995 м"
$ws.Range("D44").Value = "1 поворотов, расстояние концов 16263 м"
$ws.Range("D45").Value = "0 поворотов, расстояние концов 16263 м"
$ws.Range("D46").Value = "1 поворотов, расстояние концов 16792 м"
$ws.Range("D47").Value = "0 поворотов, расстояние концов 16757 м"
$ws.Range("D48").Value = "0 поворотов, расстояние концов 6348 м"
$ws.Range("D49").Value = "1 поворотов, расстояние концов 6348 м"
$ws.Range("D50").Value = "2 поворотов, расстояние концов 15268 м"
$ws.Range("D51").Value = "3 поворотов, расстояние концов 15268 м"
$ws.Range("D52").Value = "4 поворотов, расстояние концов 13514 м"
$ws.Range("D53").Value = "4 поворотов, расстояние концов 13514 м"
$ws.Range("D54").Value = "1 поворотов, расстояние концов 7191 м"
$ws.Range("D55").Value = "1 поворотов, расстояние концов 7194 м"
$ws.Range("D56").Value = "3 поворотов, расстояние концов 17063 м"
$ws.Range("D57").Value = "3 поворотов, расстояние концов 17134 м"
$ws.Range("D58").Value = "7 острых углов < 60°"
$ws.Range("D59").Value = "7 острых углов < 60°"
$ws.Range("D60").Value = "1 поворотов, расстояние концов 12512 м"
$ws.Range("D61").Value = "1 поворотов, расстояние концов 12512 м"
$ws.Range("D62").Value = "3 поворотов, расстояние концов 9439 м"
$ws.Range("D63").Value = "0 поворотов, расстояние концов 9430 м"
$ws.Range("D64").Value = "2 поворотов, расстояние концов 18082 м"
$ws.Range("D65").Value = "3 поворотов, расстояние концов 18144 м"
$ws.Range("D66").Value = "1 поворотов, расстояние концов 19811 м"
$ws.Range("D67").Value = "3 поворотов, расстояние концов 19801 м"
$ws.Range("D68").Value = "3 поворотов, расстояние концов 22762 м"
$ws.Range("D69").Value = "2 поворотов, расстояние концов 22764 м"
$ws.Range("D70").Value = "5 острых углов < 60°"
$ws.Range("D71").Value = "2 поворотов, расстояние концов 31144 м"
$ws.Range("D72").Value = "2 поворотов, расстояние концов 17370 м"
$ws.Range("D73").Value = "4 поворотов, расстояние концов 17223 м"
$ws.Range("D74").Value = "3 поворотов, расстояние концов 19696 м"
$ws.Range("D75").Value = "0 поворотов, расстояние концов 19696 м"
$ws.Range("D76").Value = "0 поворотов, расстояние концов 30971 м"
$ws.Range("D77").Value = "2 поворотов, расстояние концов 30971 м"
$ws.Range("D78").Value = "0 поворотов, расстояние концов 34487 м"
$ws.Range("D79").Value = "3 поворотов, расстояние концов 34503 м"
$ws.Range("D80").Value = "1 поворотов, расстояние концов 10459 м"
$ws.Range("D81").Value = "0 поворотов, расстояние концов 10459 м"
$ws.Range("D82").Value = "4 поворотов, расстояние концов 30021 м"
$ws.Range("D83").Value = "1 поворотов, расстояние концов 29872 м"
$ws.Range("D84").Value = "0 поворотов, расстояние концов 18490 м"
$ws.Range("D85").Value = "0 поворотов, расстояние концов 18338 м"
$ws.Range("D86").Value = "4 поворотов, расстояние концов 31144 м"
$ws.Range("D87").Value = "4 поворотов, расстояние концов 31144 м"
$ws.Range("D88").Value = "0 поворотов, расстояние концов 13682 м"
$ws.Range("D89").Value = "0 поворотов, расстояние концов 13666 м"
$ws.Range("D90").Value = "0 поворотов, расстояние концов 12810 м"
$ws.Range("D91").Value = "0 поворотов, расстояние концов 12786 м"
$ws.Range("D92").Value = "3 поворотов, расстояние концов 59117 м"
$ws.Range("D93").Value = "4 поворотов, расстояние концов 59112 м"
$ws.Range("D94").Value = "0 поворотов, расстояние концов 28972 м"
$ws.Range("D95").Value = "4 поворотов, расстояние концов 28972 м"
$ws.Range("D96").Value = "0 поворотов, расстояние концов 27217 м"
$ws.Range("D97").Value = "2 поворотов, расстояние концов 27285 м"
$ws.Range("D98").Value = "3 поворотов, расстояние концов 16102 м"
$ws.Range("D99").Value = "0 поворотов, расстояние концов 16077 м"
$ws.Range("D100").Value = "2 поворотов, расстояние концов 10360 м"
$ws.Range("D101").Value = "1 поворотов, расстояние концов 10500 м"
$ws.Range("D102").Value = "2 поворотов, расстояние концов 17742 м"
$ws.Range("D103").Value = "0 поворотов, расстояние концов 17746 м"
$ws.Range("D104").Value = "2 поворотов, расстояние концов 15339 м"
$ws.Range("D105").Value = "3 поворотов, расстояние концов 15339 м"
$ws.Range("D106").Value = "1 поворотов, расстояние концов 16113 м"
$ws.Range("D107").Value = "0 поворотов, расстояние концов 16119 м"
$ws.Range("D108").Value = "6 острых углов < 60°"
$ws.Range("D109").Value = "3 поворотов, расстояние концов 24522 м"
$ws.Range("D110").Value = "2 поворотов, расстояние концов 15349 м"
$ws.Range("D111").Value = "1 поворотов, расстояние концов 15339 м"
$ws.Range("D112").Value = "4 поворотов, расстояние концов 18500 м"
$ws.Range("D113").Value = "2 поворотов, расстояние концов 18352 м"
$ws.Range("D114").Value = "1 поворотов, расстояние концов 44040 м"
$ws.Range("D115").Value = "1 поворотов, расстояние концов 44045 м"
$ws.Range("D116").Value = "1 поворотов, расстояние концов 25287 м"
$ws.Range("D117").Value = "2 поворотов, расстояние концов 25287 м"
$ws.Range("D118").Value = "0 поворотов, расстояние концов 35807 м"
$ws.Range("D119").Value = "1 поворотов, расстояние концов 35824 м"
$ws.Range("D120").Value = "1 поворотов, расстояние концов 11495 м"
$ws.Range("D121").Value = "0 поворотов, расстояние концов 11459 м"
$ws.Range("D122").Value = "3 поворотов, расстояние концов 13686 м"
$ws.Range("D123").Value = "1 поворотов, расстояние концов 13892 м"
$ws.Range("D124").Value = "3 поворотов, расстояние концов 15077 м"
$ws.Range("D125").Value = "2 поворотов, расстояние концов 15077 м"
$ws.Range("D126").Value = "3 поворотов, расстояние концов 14260 м"
$ws.Range("D127").Value = "3 поворотов, расстояние концов 15114 м"
$ws.Range("D128").Value = "2 поворотов, расстояние концов 9416 м"
$ws.Range("D129").Value = "0 поворотов, расстояние концов 9345 м"
$ws.Range("D130").Value = "2 поворотов, расстояние концов 15567 м"
$ws.Range("D131").Value = "0 поворотов, расстояние концов 15567 м"
$ws.Range("D132").Value = "2 поворотов, расстояние концов 13465 м"
$ws.Range("D133").Value = "3 поворотов, расстояние концов 13505 м"
$ws.Range("D134").Value = "2 поворотов, расстояние концов 9284 м"
$ws.Range("D135").Value = "0 поворотов, расстояние концов 9351 м"
$ws.Range("D136").Value = "2 поворотов, расстояние концов 16279 м"
$ws.Range("D137").Value = "3 поворотов, расстояние концов 16282 м"
$ws.Range("D138").Value = "0 поворотов, расстояние концов 11717 м"
$ws.Range("D139").Value = "1 поворотов, расстояние концов 11717 м"
$ws.Range("D140").Value = "0 поворотов, расстояние концов 10900 м"
$ws.Range("D141").Value = "0 поворотов, расстояние концов 10900 м"
$ws.Range("D142").Value = "4 поворотов, расстояние концов 17973 м"
$ws.Range("D143").Value = "3 поворотов, расстояние концов 18089 м"
$ws.Range("D144").Value = "0 поворотов, расстояние концов 14920 м"
$ws.Range("D145").Value = "1 поворотов, расстояние концов 14917 м"
$ws.Range("D146").Value = "3 поворотов, расстояние концов 10557 м"
$ws.Range("D147").Value = "3 поворотов, расстояние концов 10557 м"
$ws.Range("D148").Value = "3 поворотов, расстояние концов 12512 м"
$ws.Range("D149").Value = "1 поворотов, расстояние концов 12512 м"
$ws.Range("D150").Value = "4 поворотов, расстояние концов 11782 м"
$ws.Range("D151").Value = "3 поворотов, расстояние концов 11689 м"
$ws.Range("D152").Value = "1 поворотов, расстояние концов 17374 м"
$ws.Range("D153").Value = "0 поворотов, расстояние концов 17374 м"
$ws.Range("D154").Value = "0 поворотов, расстояние концов 13645 м"
$ws.Range("D155").Value = "0 поворотов, расстояние концов 13584 м"
$ws.Range("D156").Value = "2 поворотов, расстояние концов 13131 м"
$ws.Range("D157").Value = "0 поворотов, расстояние концов 13008 м"
$ws.Range("D158").Value = "2 поворотов, расстояние концов 18565 м"
$ws.Range("D159").Value = "0 поворотов, расстояние концов 18565 м"
$ws.Range("D160").Value = "0 поворотов, расстояние концов 14349 м"
$ws.Range("D161").Value = "0 поворотов, расстояние концов 14349 м"
$ws.Range("D162").Value = "2 поворотов, расстояние концов 16744 м"
$ws.Range("D163").Value = "2 поворотов, расстояние концов 16746 м"
$ws.Range("D164").Value = "3 поворотов, расстояние концов 9100 м"
$ws.Range("D165").Value = "0 поворотов, расстояние концов 9067 м"
$ws.Range("D166").Value = "1 поворотов, расстояние концов 12963 м"
$ws.Range("D167").Value = "1 поворотов, расстояние концов 12942 м"
$ws.Range("D168").Value = "2 поворотов, расстояние концов 15405 м"
$ws.Range("D169").Value = "2 поворотов, расстояние концов 15376 м"
$ws.Range("D170").Value = "3 поворотов, расстояние концов 10279 м"
$ws.Range("D171").Value = "1 поворотов, расстояние концов 10279 м"
$ws.Range("D172").Value = "0 поворотов, расстояние концов 14866 м"
$ws.Range("D173").Value = "2 поворотов, расстояние концов 14866 м"
$ws.Range("D174").Value = "1 поворотов, расстояние концов 14513 м"
$ws.Range("D175").Value = "1 поворотов, расстояние концов 14518 м"
$ws.Range("D177").Value = "расстояние между концами 29 м"
$ws.Range("D179").Value = "расстояние между концами 29 м"
$ws.Range("D180").Value = "0 поворотов, расстояние концов 10088 м"
$ws.Range("D181").Value = "0 поворотов, расстояние концов 10088 м"
$ws.Range("D182").Value = "0 поворотов, расстояние концов 10125 м"
$ws.Range("D183").Value = "1 поворотов, расстояние концов 10125 м"
$ws.Range("D184").Value = "3 поворотов, расстояние концов 13019 м"
$ws.Range("D185").Value = "2 поворотов, расстояние концов 12925 м"
$ws.Range("D186").Value = "6 острых углов < 60°"
$ws.Range("D187").Value = "0 поворотов, расстояние концов 28031 м"
$ws.Range("D188").Value = "0 поворотов, расстояние концов 13584 м"
$ws.Range("D189").Value = "3 поворотов, расстояние концов 13584 м"
$ws.Range("D190").Value = "1 поворотов, расстояние концов 6789 м"
$ws.Range("D191").Value = "0 поворотов, расстояние концов 6789 м"
$ws.Range("D192").Value = "1 поворотов, расстояние концов 3824 м"
$ws.Range("D193").Value = "5 острых углов < 60°"
$ws.Range("D194").Value = "3 поворотов, расстояние концов 8933 м"
$ws.Range("D195").Value = "0 поворотов, расстояние концов 8924 м"
$ws.Range("D196").Value = "1 поворотов, расстояние концов 22068 м"
$ws.Range("D197").Value = "0 поворотов, расстояние концов 22049 м"
$ws.Range("D198").Value = "1 поворотов, расстояние концов 8502 м"
$ws.Range("D199").Value = "0 поворотов, расстояние концов 8502 м"
$ws.Range("D200").Value = "2 поворотов, расстояние концов 12404 м"
$ws.Range("D201").Value = "3 поворотов, расстояние концов 12390 м"
$ws.Range("D202").Value = "2 поворотов, расстояние концов 19873 м"
$ws.Range("D203").Value = "1 поворотов, расстояние концов 19873 м"
$ws.Range("D204").Value = "1 поворотов, расстояние концов 9087 м"
$ws.Range("D205").Value = "2 поворотов, расстояние концов 8960 м"
$ws.Range("D206").Value = "1 поворотов, расстояние концов 10435 м"
$ws.Range("D207").Value = "2 поворотов, расстояние концов 10316 м"
$ws.Range("D208").Value = "0 поворотов, расстояние концов 11787 м"
$ws.Range("D209").Value = "0 поворотов, расстояние концов 11680 м"
$ws.Range("D210").Value = "0 поворотов, расстояние концов 11670 м"
$ws.Range("D211").Value = "1 поворотов, расстояние концов 11713 м"
$ws.Range("D212").Value = "3 поворотов, расстояние концов 16979 м"
$ws.Range("D213").Value = "3 поворотов, расстояние концов 16979 м"
$ws.Range("D214").Value = "3 поворотов, расстояние концов 11635 м"
$ws.Range("D215").Value = "5 острых углов < 60°"
$ws.Range("D216").Value = "4 поворотов, расстояние концов 13620 м"
$ws.Range("D217").Value = "0 поворотов, расстояние концов 13514 м"
$ws.Range("D218").Value = "3 поворотов, расстояние концов 14674 м"
$ws.Range("D219").Value = "1 поворотов, расстояние концов 14674 м"
$ws.Range("D220").Value = "4 поворотов, расстояние концов 14923 м"
$ws.Range("D221").Value = "2 поворотов, расстояние концов 14920 м"
$ws.Range("D222").Value = "0 поворотов, расстояние концов 11855 м"
$ws.Range("D223").Value = "0 поворотов, расстояние концов 11855 м"
$ws.Range("D224").Value = "3 поворотов, расстояние концов 10100 м"
$ws.Range("D225").Value = "2 поворотов, расстояние концов 10208 м"
$ws.Range("D226").Value = "2 поворотов, расстояние концов 5153 м"
$ws.Range("D227").Value = "4 поворотов, расстояние концов 5146 м"
$ws.Range("D228").Value = "0 поворотов, расстояние концов 16174 м"
$ws.Range("D229").Value = "0 поворотов, расстояние концов 16174 м"
